$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Table cell left margins: 123 dxa (6.15pt) -> 128 dxa (6.4pt)
#    Applies to both tables in the body.
# ---------------------------------------------------------------------
foreach ($tbl in $d.Tables) {
    $tbl.LeftPadding = 6.4
}

# ---------------------------------------------------------------------
# 2) Merge-field rename: contact.* -> recipient.*
# ---------------------------------------------------------------------

# [contact.postal_address;strconv=no] -> [recipient.postal_address;strconv=no]
$d.Content.Find.Execute(
    "[contact.postal_address;strconv=no]", $false, $true, $false, $false, $false,
    $false, 1, $false, "[recipient.postal_address;strconv=no]", 2) | Out-Null

# [attachments.chrono] -> [attachment.chrono]
$d.Content.Find.Execute(
    "[attachments.chrono]", $false, $true, $false, $false, $false,
    $false, 1, $false, "[attachment.chrono]", 2) | Out-Null

# [contact.contact_title] [contact.contact_lastname], -> [recipient.civility] [recipient.lastname],
$d.Content.Find.Execute(
    "[contact.contact_title] [contact.contact_lastname],", $false, $true, $false, $false, $false,
    $false, 1, $false, "[recipient.civility] [recipient.lastname],", 2) | Out-Null

# Veuillez agreer, [contact.contact_title], ... -> Veuillez agreer, [recipient.civility], ...
$d.Content.Find.Execute(
    "Veuillez agréer, [contact.contact_title], l’expression de nos salutations distinguées.", $false, $true, $false, $false, $false,
    $false, 1, $false, "Veuillez agréer, [recipient.civility], l’expression de nos salutations distinguées.", 2) | Out-Null

# ---------------------------------------------------------------------
# 3) Update the cached TIME field result: 31/10/2019 -> 09/12/2019
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "31/10/2019", $false, $true, $false, $false, $false,
    $false, 1, $false, "09/12/2019", 2) | Out-Null

# ---------------------------------------------------------------------
# 4) Header separator line drawing: resize slightly
#    wp:extent/a:ext cx,cy grow by a few EMU.
# ---------------------------------------------------------------------
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
foreach ($shp in $hdr.Shapes) {
    if ($shp.Name -eq "Image1" -and [Math]::Round($shp.Width, 1) -eq 543.8) {
        $shp.Width = [single](6906895.0 / 12700.0)
        $shp.Height = [single](3810.0 / 12700.0)
    }
}
